$d = $word.ActiveDocument
$wdNS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ------------------------------------------------------------------
# Locate the paragraph that ends "...как бы это ни было печально."
# (the "_GoBack" bookmark currently sits right at the end of it).
# ------------------------------------------------------------------
$probe = $d.Content
$found = $probe.Find.Execute("как бы это ни было печально.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraIndex = $probe.Paragraphs(1).Index

$targetPara = $d.Paragraphs($paraIndex)
$insPos = $targetPara.Range.End
$insertionPoint = $d.Range($insPos, $insPos)

# ------------------------------------------------------------------
# Insert the four new paragraphs after it. A trailing "anchor"
# paragraph is appended to the XML batch purely so the real content
# keeps its own formatting (InsertXML merges the *last* fragment into
# whatever paragraph follows the insertion point); the anchor's stray
# text is then stripped away, leaving the original trailing empty
# paragraph mark (and its formatting) untouched, ready to carry the
# relocated bookmark.
# ------------------------------------------------------------------
$newParasXml = @"
<w:p $wdNS><w:pPr><w:pStyle w:val="a0"/></w:pPr><w:r><w:t>На следующей странице перейдём к обзору одного из этапов теста. Исходный код остальных тестов можете найти по ссылке, данной в списке литературы.</w:t></w:r></w:p><w:p $wdNS><w:pPr><w:pStyle w:val="a0"/></w:pPr></w:p><w:p $wdNS><w:pPr><w:pStyle w:val="a0"/><w:jc w:val="center"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Этап 2. Атака на основе открытых текстов и соответствующих им шифротекстов.</w:t></w:r></w:p><w:p $wdNS><w:pPr><w:pStyle w:val="a0"/></w:pPr><w:r><w:t>На самом деле же пара нам понадобится лишь одна. Из-за всё того же маленького диапазона возможных значений ключа. Рассмотрим алгоритм взлома подробнее.</w:t></w:r></w:p><w:p $wdNS><w:r><w:t>ANCHOR_PARAGRAPH_TO_REMOVE</w:t></w:r></w:p>
"@

$insertionPoint.InsertXML($newParasXml)

# ------------------------------------------------------------------
# Strip the throw-away anchor paragraph's text back out, leaving the
# formerly-trailing empty paragraph exactly as it was.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$anchorTextRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$anchorTextRange.Delete()

# ------------------------------------------------------------------
# Move the "_GoBack" bookmark onto that same (now) final paragraph.
# ------------------------------------------------------------------
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$bmRange = $finalPara.Range
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
